$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "No Emulator Data"
$ws.Range("C13").Value = "12"
$ws.Range("D13").Value = "2.5"
$ws.Range("E13").Value = "10.50"
$ws.Range("G13").Value = "PayNow"
$ws.Range("H13").Value = "en_US"
$ws.Range("W13").Value = "udf data 4"
$ws.Range("X13").Value = "udf data 5"
$ws.Range("AB13").Value = "udf data 9"
$ws.Range("AC13").Value = "udf data 10"

$ws.Range("S11").Copy()
$ws.Range("S13").PasteSpecial(-4122)

$ws.Range("C13").Select()
